$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trebuchet-V1-BOM")

# --- Update Designation for TCMT4100 -> TLP293-4 (row 8) ---
$ws.Range("E8").Value = "TLP293-4"

# --- Update the terminal-block comment text (row 14) to mention the LCSC count note ---
$ws.Range("G14").Value = "Make sure to use terminals with 5mm spacing not 5.08mm! Does not need to be 10 at one piece. Combining multiple to add to 10 will be fine. (LSCS number is only 2 so 5 need to be connected together)"

# --- Add new "LCSC Number" column (J) ---
$ws.Range("J1").Value = "LCSC Number"

$ws.Range("J2").Value = "C21190"
$ws.Range("J3").Value = "C129022"
$ws.Range("J4").Value = "C22962"
$ws.Range("J5").Value = "C14663"
$ws.Range("J6").Value = "C25804"
$ws.Range("J7").Value = "C845537"
$ws.Range("J8").Value = "C112623"
$ws.Range("J10").Value = "C23212"
$ws.Range("J11").Value = "C2897391"
$ws.Range("J12").Value = "C2894966"
$ws.Range("J13").Value = "C506653"
$ws.Range("J14").Value = "C8404"
